# Version 1.4 in progress
# The mail functionality as well as the side bar were unified with 2 working charts
#
# Update the "Provincias"/"Celular"/"Correo" data for the associates table and
# wire up the missing mailto: hyperlink for Cesar's e-mail address.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell contents -------------------------------------------------
# Order matters: new shared-string entries are appended in first-seen order,
# so we touch the cells in the same sequence the workbook author did.

# Row 4 - Cesar Jimenez: fill in the previously-empty e-mail cell and his
# phone number.
$ws.Range("B4").Value = "cesarjjxd@gmail.com"
$ws.Range("C4").Value = "8529-6827"

# Row 3 - Fernanda Murillo: new phone number.
$ws.Range("C3").Value = "8598-6048"

# Row 4 - Cesar Jimenez: provinces he can cover.
$ws.Range("D4").Value = "San José,Heredia,Cartago"

# Row 2 - Maynor Martinez: extended province coverage.
$ws.Range("D2").Value = "San José,Heredia,Cartago,Alajuela,Puntarenas"

# Row 3 - Fernanda Murillo: same extended province coverage.
$ws.Range("D3").Value = "San José,Heredia,Cartago,Alajuela,Puntarenas"

# --- Hyperlink for the newly added e-mail address --------------------------
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:cesarjjxd@gmail.com")

# --- Column D is now much wider to fit the longer province lists -----------
$ws.Range("D1").EntireColumn.ColumnWidth = 47.65

# --- Leave the selection on D2, like in the saved workbook -----------------
$ws.Range("D2").Select()
